# Updated cryptos list on Tue Apr 23 20:39:15 UTC 2024 with GitHub Actions
# Refresh Price (D) / Volume(1h) (E) columns, and Bittensor/OKB swap ranks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "66.212.08";  E = "  -0.45%  " },
    @{ Row = 3;  D = "3.202.01";   E = "  +0.28%  " },
    @{ Row = 4;  D = $null;        E = "  -0.12%  " },
    @{ Row = 5;  D = "608.02";     E = "  +1.92%  " },
    @{ Row = 6;  D = "156.64";     E = "  +1.39%  " },
    @{ Row = 7;  D = $null;        E = "  -0.01%  " },
    @{ Row = 8;  D = "3.200.99";   E = "  +0.34%  " },
    @{ Row = 9;  D = "0.554";      E = "  -1.57%  " },
    @{ Row = 10; D = "0.161";      E = "  -0.23%  " },
    @{ Row = 11; D = "5.68";       E = "  -4.13%  " },
    @{ Row = 12; D = "0.505";      E = "  -2.87%  " },
    @{ Row = 13; D = "0.0000270";  E = "  +0.76%  " },
    @{ Row = 14; D = "38.48";      E = "  -2.08%  " },
    @{ Row = 15; D = "3.726.59";   E = "  +0.33%  " },
    @{ Row = 16; D = "66.331.74";  E = "  -0.27%  " },
    @{ Row = 17; D = "7.36";       E = "  -1.90%  " },
    @{ Row = 18; D = "3.202.40";   E = "  +0.27%  " },
    @{ Row = 19; D = $null;        E = "  +1.47%  " },
    @{ Row = 20; D = "508.78";     E = "  -1.62%  " },
    @{ Row = 21; D = "15.34";      E = "  -0.60%  " },
    @{ Row = 22; D = "0.734";      E = "  -0.69%  " },
    @{ Row = 23; D = "8.02";       E = "  -1.05%  " },
    @{ Row = 24; D = "14.66";      E = "  -2.07%  " },
    @{ Row = 25; D = "85.14";      E = "  -1.05%  " },
    @{ Row = 26; D = $null;        E = "  -0.06%  " },
    @{ Row = 27; D = "3.00";       E = "  -0.17%  " },
    @{ Row = 28; D = "9.08";       E = "  -2.15%  " },
    @{ Row = 29; D = "2.36";       E = "  +0.10%  " },
    @{ Row = 30; D = "0.128";      E = "  +41.36%  " },
    @{ Row = 31; D = "2.94";       E = "  +0.31%  " },
    @{ Row = 32; D = "7.00";       E = "  -1.26%  " },
    @{ Row = 33; D = "28.29";      E = "  -0.40%  " },
    @{ Row = 34; D = $null;        E = "  +0.00%  " },
    @{ Row = 35; D = $null;        E = "  -4.68%  " },
    @{ Row = 36; D = "6.51";       E = "  -0.59%  " },
    @{ Row = 40; D = $null;        E = "  +2.44%  " },
    @{ Row = 41; D = "0.0422";     E = "  -0.99%  " },
    @{ Row = 42; D = "3.06";       E = "  +5.75%  " },
    @{ Row = 43; D = "8.76";       E = "  -1.88%  " },
    @{ Row = 44; D = "0.297";      E = "  -1.47%  " },
    @{ Row = 45; D = "2.910.45";   E = "  +0.05%  " },
    @{ Row = 46; D = "2.44";       E = "  -0.46%  " },
    @{ Row = 47; D = "28.33";      E = "  -1.91%  " },
    @{ Row = 48; D = "2.41";       E = "  +2.02%  " },
    @{ Row = 49; D = $null;        E = "  -0.01%  " },
    @{ Row = 50; D = "0.117";      E = "  -1.03%  " },
    @{ Row = 51; D = "122.53";     E = "  +0.05%  " }
)

# Rows whose new Price text is numeric-looking ("608.02", "0.0000270", ...).
# Mark column D as Text *before* writing so Excel stores the literal digits
# instead of silently coercing them to a floating point number.
$textRows = @(5,6,9,10,11,12,13,14,17,20,21,22,23,24,25,27,28,29,30,31,32,33,36,41,42,43,44,46,47,48,50,51)
foreach ($r in $textRows) {
    $ws.Range("D$r").NumberFormat = "@"
}
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $ws.Range("D$r").Value = $u.D
    }
    $ws.Range("E$r").Value = $u.E
}

# D39 keeps a subscript-3 glyph between the leading zeros: 0.0\u20830777
$sub3 = [char]0x2083
$ws.Range("D39").Value = [string]::Concat("0.0", $sub3, "0777")
$ws.Range("E39").Value = "  +15.85%  "

# Rows 37 and 38 swap rank positions: Bittensor overtakes OKB.
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").Value = "502.10"
$ws.Range("E37").Value = "  -1.73%  "

$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").Value = "55.53"
$ws.Range("E38").Value = "  +0.91%  "
